$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 68332.336
$ws.Range("I28").Value = 68332.336
$ws.Range("K28").Value = 68332.336
$ws.Range("M28").Value = -67847.336
# Row 33
$ws.Range("H33").Value = 410.25
$ws.Range("I33").Value = 104.26667
$ws.Range("J33").Value = 5000
$ws.Range("K33").Value = 104.26667
$ws.Range("L33").Value = 5000
$ws.Range("M33").Value = 124.73333
$ws.Range("N33").Value = -5458
# Row 40
$ws.Range("H40").Value = 7204.1875
$ws.Range("J40").Value = 8086.5654
$ws.Range("L40").Value = 8086.5654
$ws.Range("N40").Value = -8436.565399999999
# Row 43
$ws.Range("H43").Value = 3333
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 3999.5
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 3999.5
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -4137.5
# Row 64
$ws.Range("H64").Value = 18004100
$ws.Range("J64").Value = 5166.6665
$ws.Range("L64").Value = 5166.6665
$ws.Range("N64").Value = -5662.6665
# Row 67
$ws.Range("H67").Value = 18004100
$ws.Range("J67").Value = 5166.6665
$ws.Range("L67").Value = 5166.6665
$ws.Range("N67").Value = -6882.6665
# Row 80
$ws.Range("H80").Value = 8266.625
$ws.Range("I80").Value = 560
$ws.Range("K80").Value = 1680
$ws.Range("M80").Value = -682
# Row 83
$ws.Range("H83").Value = 8266.625
$ws.Range("I83").Value = 560
$ws.Range("K83").Value = 5040
$ws.Range("M83").Value = -48
# Row 95
$ws.Range("H95").Value = 51880
$ws.Range("J95").Value = 51880
$ws.Range("L95").Value = 51880
$ws.Range("N95").Value = -57372
# Row 112
$ws.Range("H112").Value = 4350
$ws.Range("J112").Value = 4362.5
$ws.Range("L112").Value = 13087.5
$ws.Range("N112").Value = -15303.5
# Row 141
$ws.Range("H141").Value = 14260.667
$ws.Range("I141").Value = 9114.1
$ws.Range("K141").Value = 27342.3
$ws.Range("M141").Value = -22162.3

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("N76").Value = 0
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").Value = 0

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 120
$ws.Range("H120").Value = 84997
$ws.Range("J120").Value = 84997
$ws.Range("L120").Value = 84997
$ws.Range("N120").Value = -94673

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 355.75
$ws.Range("J7").Value = 333.33334
$ws.Range("L7").Value = 333.33334
$ws.Range("N7").Value = -559.33334
# Row 50
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31250
# Row 51
$ws.Range("H51").Value = 17500.5
# Row 60
$ws.Range("H60").Value = 10166.667
$ws.Range("J60").Value = 22500
$ws.Range("L60").Value = 22500
$ws.Range("N60").Value = -23522
# Row 61
$ws.Range("H61").Value = 17500.5
# Row 62
$ws.Range("H62").Value = 11994.75
$ws.Range("I62").Value = 11994.75
$ws.Range("K62").Value = 11994.75
$ws.Range("M62").Value = -11370.75
# Row 65
$ws.Range("H65").Value = 11994.75
$ws.Range("I65").Value = 11994.75
$ws.Range("K65").Value = 59973.75
$ws.Range("M65").Value = -56853.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 107020.8
$ws.Range("I12").Value = 214025.6
$ws.Range("J12").Value = 16
$ws.Range("K12").Value = 642076.8
$ws.Range("L12").Value = 48
$ws.Range("M12").Value = -641903.8
$ws.Range("N12").Value = -394
# Row 39
$ws.Range("H39").Value = 2960.5
$ws.Range("J39").Value = 3287.4285
$ws.Range("L39").Value = 9862.2855
$ws.Range("N39").Value = -10450.2855
# Row 129
$ws.Range("H129").Value = 2751.125
$ws.Range("I129").Value = 1713.1428
$ws.Range("J129").Value = 3558.4443
$ws.Range("K129").Value = 5139.428400000001
$ws.Range("L129").Value = 10675.3329
$ws.Range("M129").Value = -139.4284000000007
$ws.Range("N129").Value = -20675.3329
# Row 140
$ws.Range("H140").Value = 1095.5
$ws.Range("I140").Value = 708
$ws.Range("J140").Value = 3033
$ws.Range("K140").Value = 2124
$ws.Range("L140").Value = 9099
$ws.Range("M140").Value = 3056
$ws.Range("N140").Value = -19459

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 357814.16
$ws.Range("J3").Value = 674.75
$ws.Range("L3").Value = 674.75
$ws.Range("N3").Value = -906.75
# Row 55
$ws.Range("H55").Value = 16999
$ws.Range("J55").Value = 24998.334
$ws.Range("L55").Value = 24998.334
$ws.Range("N55").Value = -25652.334
# Row 105
$ws.Range("H105").Value = 20890.334
$ws.Range("J105").Value = 20890.334
$ws.Range("L105").Value = 20890.334
$ws.Range("N105").Value = -27878.334

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3004
$ws.Range("I16").Value = 1755
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 1755
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -1585
$ws.Range("N16").Value = -8340
# Row 17
$ws.Range("H17").Value = 2127
$ws.Range("I17").Value = 2127
$ws.Range("K17").Value = 2127
$ws.Range("M17").Value = -1957
# Row 22
$ws.Range("H22").Value = 2691.5652
$ws.Range("I22").Value = 2469
$ws.Range("J22").Value = 2810.2666
$ws.Range("K22").Value = 2469
$ws.Range("L22").Value = 2810.2666
$ws.Range("M22").Value = -2174
$ws.Range("N22").Value = -3400.2666
# Row 27
$ws.Range("H27").Value = 2691.5652
$ws.Range("I27").Value = 2469
$ws.Range("J27").Value = 2810.2666
$ws.Range("K27").Value = 2469
$ws.Range("L27").Value = 2810.2666
$ws.Range("M27").Value = -2362
$ws.Range("N27").Value = -3024.2666
# Row 46
$ws.Range("H46").Value = 3782.7273
$ws.Range("I46").Value = 3357.0715
$ws.Range("K46").Value = 3357.0715
$ws.Range("M46").Value = -3169.0715
# Row 68
$ws.Range("H68").Value = 2859.2666
$ws.Range("J68").Value = 1199.5
$ws.Range("L68").Value = 1199.5
$ws.Range("N68").Value = -2697.5
# Row 71
$ws.Range("H71").Value = 2859.2666
$ws.Range("J71").Value = 1199.5
$ws.Range("L71").Value = 5997.5
$ws.Range("N71").Value = -13485.5
# Row 81
$ws.Range("H81").Value = 80000
$ws.Range("I81").Value = 80000
$ws.Range("K81").Value = 80000
$ws.Range("M81").Value = -79002
# Row 84
$ws.Range("H84").Value = 80000
$ws.Range("I84").Value = 80000
$ws.Range("K84").Value = 240000
$ws.Range("M84").Value = -235008

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 99995
$ws.Range("J16").Value = 99995
$ws.Range("L16").Value = 99995
$ws.Range("N16").Value = -100579
# Row 58
$ws.Range("H58").Value = 60000
$ws.Range("I58").Value = 60000
$ws.Range("K58").Value = 60000
$ws.Range("M58").Value = -59692
# Row 64
$ws.Range("H64").Value = 79999
$ws.Range("I64").Value = 79999
$ws.Range("K64").Value = 79999
$ws.Range("M64").Value = -79751
# Row 67
$ws.Range("H67").Value = 79999
$ws.Range("I67").Value = 79999
$ws.Range("K67").Value = 79999
$ws.Range("M67").Value = -79141
# Row 75
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
# Row 78
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
